# Swap columns C (codeforiati:group-name) and D (codeforiati:group-code)
# for every row in the sheet, including the header row, to match the
# updated column ordering in the source codelist (group-code before
# group-name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cVal = $cCell.Value()
    $dVal = $dCell.Value()

    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
